$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 2
$ws.Range("B6").Value = 5
$ws.Range("B7").Value = 10
$ws.Range("B8").Value = 28

$ws.Range("N2").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("N6").Value = 1
$ws.Range("N7").Value = 5
$ws.Range("N8").Value = 29

$ws.Range("M16").Select()
